# Update the "想去人数" (number of people wanting to go) counts in the
# "展览" and "全部类型" sheets to reflect newly scraped values.

$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F, applied identically in both sheets
# (the rows correspond to the same events, just listed at different row
# offsets because "全部类型" interleaves rows from other sheets as well).

$wsExhibition = $wb.Worksheets.Item("展览")
$updatesExhibition = @{
    2  = 153
    6  = 1324
    12 = 216
    15 = 473
    16 = 93
    18 = 497
    19 = 291
    20 = 414
    21 = 106
    22 = 226
    26 = 450
    27 = 302
}
foreach ($row in $updatesExhibition.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $updatesExhibition[$row]
}

$wsAll = $wb.Worksheets.Item("全部类型")
$updatesAll = @{
    4  = 153
    8  = 1324
    19 = 216
    22 = 473
    23 = 93
    25 = 497
    28 = 291
    29 = 414
    31 = 106
    33 = 226
    41 = 450
    42 = 302
}
foreach ($row in $updatesAll.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $updatesAll[$row]
}
